# Final working UiPath project
# The "Test Results" sheet tracked pass/fail status for CRUD tests.
# The "Delete Test Passed" column (E) for rows 15-23 was flipped from
# TRUE to FALSE (these were previously marked as passed, now corrected
# to reflect that the delete tests did not pass).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Results")

$ws.Range("E15:E23").Value = $false
